$p = $ppt.ActivePresentation

# --- Slide 11: replace the free-floating "SEKIAN DAN TERIMAKASIH" textbox
#     with the layout's centre-title placeholder, now reading "THANK YOU".
$s11 = $p.Slides.Item(11)

# Remove the old plain textbox entirely.
$s11.Shapes.Item(1).Delete()

# A throwaway textbox to "consume" shape id 2, so the placeholder we paste
# next lands on id 3 (matching the target deck).
$dummy = $s11.Shapes.AddTextbox(1, 0, 0, 10, 10)

# Borrow a clean ctrTitle placeholder (no stray hasCustomPrompt / lstStyle
# overrides) from slide 1 and drop it onto slide 11.
$s1 = $p.Slides.Item(1)
$titleSource = $s1.Shapes.Item(1)
$titleSource.Copy()
$pasted = $s11.Shapes.Paste()
$title = $pasted.Item(1)

$dummy.Delete()

$title.Name = "Title 19"

$title.Left = 144.93914794921875
$title.Top = 166.36166381835938
$title.Width = 535.2952880859375
$title.Height = 131.84803771972656

$title.TextFrame.TextRange.Delete()
$title.TextFrame.TextRange.Text = "THANK YOU"
